$d = $word.ActiveDocument

function Get-ParaRange($idx) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($i -eq $idx) { return $p.Range }
    }
}

# ---------------------------------------------------------------------
# 1. Bold the "Q"/"1" (but not the ":") in the "Q1:" heading (paragraph 1)
# ---------------------------------------------------------------------
$d.Range(0, 2).Bold = 1

# ---------------------------------------------------------------------
# 2. Bold the "A1" (but not the ":") in the "A1:" line (paragraph 6)
# ---------------------------------------------------------------------
$d.Range(332, 334).Bold = 1

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the Q2-list paragraph
#    up onto the blank paragraph that follows "A1:" (paragraph 7)
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$p7 = Get-ParaRange 7
$p7.Bookmarks.Add("_GoBack") | Out-Null

# ---------------------------------------------------------------------
# 4. Strip the stray "rFonts hint=eastAsia" paragraph-mark formatting from
#    the blank paragraph right before "Q2:" (paragraph 9), turning it into
#    a bare, unformatted paragraph.
# ---------------------------------------------------------------------
$p9 = Get-ParaRange 9
$p9.InsertXML("<dummy/>")

# ---------------------------------------------------------------------
# 5. Bold the "Q"/"2" (but not the ":") in the "Q2:" heading (paragraph 10)
# ---------------------------------------------------------------------
$d.Range(339, 341).Bold = 1

# ---------------------------------------------------------------------
# 6. Bold the "A2" (but not the ":") in the "A2:" line (paragraph 15)
# ---------------------------------------------------------------------
$d.Range(736, 738).Bold = 1

# ---------------------------------------------------------------------
# 7. Strip the stray "rFonts hint=eastAsia" paragraph-mark formatting from
#    the trailing blank paragraph (paragraph 16), turning it into a bare,
#    unformatted paragraph.
# ---------------------------------------------------------------------
$p16 = Get-ParaRange 16
$p16.InsertXML("<dummy/>")
